$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: insert a new column before D (old Dec_2025) to make room for Jan_2026,
# shifting Dec_2025/Nov_2025/Oct_2025/MoM/QoQ one column to the right (D->E,E->F,F->G,G->H,H->I)
$ws.Columns.Item(4).Insert()

# Step 2: drop the old Oct_2025 column, now shifted to column G, so MoM/QoQ shift back to G/H
$ws.Columns.Item(7).Delete()

# Step 3: the source table lost one holding (27 data rows -> 25 data rows); delete the last row
# (current row 27) so the sheet shrinks from 27 rows to 26 rows (header + 25 data rows).
$ws.Rows.Item(27).Delete()

# Step 4: header row
$ws.Cells.Item(1,1).Value = "ISIN"
$ws.Cells.Item(1,2).Value = "Stock Name"
$ws.Cells.Item(1,3).Value = "Mutual Fund"
$ws.Cells.Item(1,4).Value = "Jan_2026"
$ws.Cells.Item(1,5).Value = "Dec_2025"
$ws.Cells.Item(1,6).Value = "Nov_2025"
$ws.Cells.Item(1,7).Value = "MoM"
$ws.Cells.Item(1,8).Value = "QoQ"

# Step 5: refreshed holdings data from the quant engine (re-sorted by Jan_2026 weight,
# two holdings dropped, two new holdings added)
$ws.Cells.Item(2,1).Value = "INE406A01037"
$ws.Cells.Item(2,2).Value = "Aurobindo Pharma Limited"
$ws.Cells.Item(2,3).Value = "quant Mid Cap Fund"
$ws.Cells.Item(2,4).Value = 9.897539
$ws.Cells.Item(2,5).Value = 8.763819
$ws.Cells.Item(2,6).Value = 8.767438
$ws.Cells.Item(2,7).Value = 1.13372
$ws.Cells.Item(2,8).Value = 1.130101
$ws.Cells.Item(3,1).Value = "INE151A01013"
$ws.Cells.Item(3,2).Value = "Tata Communications Limited"
$ws.Cells.Item(3,3).Value = "quant Mid Cap Fund"
$ws.Cells.Item(3,4).Value = 8.547462
$ws.Cells.Item(3,5).Value = 8.979951
$ws.Cells.Item(3,6).Value = 8.615686
$ws.Cells.Item(3,7).Value = -0.4324890000000003
$ws.Cells.Item(3,8).Value = -0.06822400000000073
$ws.Cells.Item(4,1).Value = "INE821I01022"
$ws.Cells.Item(4,2).Value = "IRB Infrastructure Developers Limited"
$ws.Cells.Item(4,3).Value = "quant Mid Cap Fund"
$ws.Cells.Item(4,4).Value = 7.117608
$ws.Cells.Item(4,5).Value = 6.668885
$ws.Cells.Item(4,6).Value = 6.565554
$ws.Cells.Item(4,7).Value = 0.4487229999999993
$ws.Cells.Item(4,8).Value = 0.552054
$ws.Cells.Item(5,1).Value = "INE417T01026"
$ws.Cells.Item(5,2).Value = "PB Fintech Limited"
$ws.Cells.Item(5,3).Value = "quant Mid Cap Fund"
$ws.Cells.Item(5,4).Value = 6.003223
$ws.Cells.Item(5,5).Value = 5.987753
$ws.Cells.Item(5,6).Value = 2.143951
$ws.Cells.Item(5,7).Value = 0.01547000000000054
$ws.Cells.Item(5,8).Value = 3.859272
$ws.Cells.Item(6,1).Value = "INE0BS701011"
$ws.Cells.Item(6,2).Value = "Premier Energies Limited"
$ws.Cells.Item(6,3).Value = "quant Mid Cap Fund"
$ws.Cells.Item(6,4).Value = 5.903756
$ws.Cells.Item(6,5).Value = 6.242897
$ws.Cells.Item(6,6).Value = 6.474772
$ws.Cells.Item(6,7).Value = -0.3391410000000006
$ws.Cells.Item(6,8).Value = -0.5710160000000002
$ws.Cells.Item(7,1).Value = "INE281B01032"
$ws.Cells.Item(7,2).Value = "Lloyds Metals And Energy Limited"
$ws.Cells.Item(7,3).Value = "quant Mid Cap Fund"
$ws.Cells.Item(7,4).Value = 5.563413
$ws.Cells.Item(7,5).Value = 6.01995
$ws.Cells.Item(7,6).Value = 5.364377
$ws.Cells.Item(7,7).Value = -0.456537
$ws.Cells.Item(7,8).Value = 0.1990359999999995
$ws.Cells.Item(8,1).Value = "INE0CZ201020"
$ws.Cells.Item(8,2).Value = "ANTHEM BIOSCIENCES LIMITED"
$ws.Cells.Item(8,3).Value = "quant Mid Cap Fund"
$ws.Cells.Item(8,4).Value = 4.177577
$ws.Cells.Item(8,5).Value = 4.105208
$ws.Cells.Item(8,6).Value = 2.431707
$ws.Cells.Item(8,7).Value = 0.07236900000000013
$ws.Cells.Item(8,8).Value = 1.74587
$ws.Cells.Item(9,1).Value = "INE473A01011"
$ws.Cells.Item(9,2).Value = "Linde India Ltd."
$ws.Cells.Item(9,3).Value = "quant Mid Cap Fund"
$ws.Cells.Item(9,4).Value = 3.948508
$ws.Cells.Item(9,5).Value = 3.538621
$ws.Cells.Item(9,6).Value = 3.41999
$ws.Cells.Item(9,7).Value = 0.4098869999999999
$ws.Cells.Item(9,8).Value = 0.528518
$ws.Cells.Item(10,1).Value = "INE018E01016"
$ws.Cells.Item(10,2).Value = "SBI Cards & Payment Services Ltd"
$ws.Cells.Item(10,3).Value = "quant Mid Cap Fund"
$ws.Cells.Item(10,4).Value = 3.113632
$ws.Cells.Item(10,5).Value = 2.845461
$ws.Cells.Item(10,6).Value = 1.529238
$ws.Cells.Item(10,7).Value = 0.2681710000000002
$ws.Cells.Item(10,8).Value = 1.584394
$ws.Cells.Item(11,1).Value = "INE880J01026"
$ws.Cells.Item(11,2).Value = "JSW Infrastructure Limited"
$ws.Cells.Item(11,3).Value = "quant Mid Cap Fund"
$ws.Cells.Item(11,4).Value = 2.346329
$ws.Cells.Item(11,5).Value = 2.336732
$ws.Cells.Item(11,6).Value = 2.140035
$ws.Cells.Item(11,7).Value = 0.009596999999999856
$ws.Cells.Item(11,8).Value = 0.2062939999999998
$ws.Cells.Item(12,1).Value = "INE002A01018"
$ws.Cells.Item(12,2).Value = "Reliance Industries Limited"
$ws.Cells.Item(12,3).Value = "quant Mid Cap Fund"
$ws.Cells.Item(12,4).Value = 2.030154
$ws.Cells.Item(12,5).Value = 6.012765
$ws.Cells.Item(12,6).Value = 5.790248
$ws.Cells.Item(12,7).Value = -3.982611
$ws.Cells.Item(12,8).Value = -3.760094
$ws.Cells.Item(13,1).Value = "INE484J01027"
$ws.Cells.Item(13,2).Value = "Godrej Properties Limited"
$ws.Cells.Item(13,3).Value = "quant Mid Cap Fund"
$ws.Cells.Item(13,4).Value = 1.593499
$ws.Cells.Item(13,5).Value = 1.831049
$ws.Cells.Item(13,6).Value = 0.880814
$ws.Cells.Item(13,7).Value = -0.2375499999999999
$ws.Cells.Item(13,8).Value = 0.712685
$ws.Cells.Item(14,1).Value = "INE042A01014"
$ws.Cells.Item(14,2).Value = "Escorts Kubota Limited"
$ws.Cells.Item(14,3).Value = "quant Mid Cap Fund"
$ws.Cells.Item(14,4).Value = 1.570087
$ws.Cells.Item(14,5).Value = 1.561581
$ws.Cells.Item(14,6).Value = 1.546102
$ws.Cells.Item(14,7).Value = 0.008505999999999903
$ws.Cells.Item(14,8).Value = 0.02398499999999992
$ws.Cells.Item(15,1).Value = "INE686F01025"
$ws.Cells.Item(15,2).Value = "UNITED BREWERIES LIMITED"
$ws.Cells.Item(15,3).Value = "quant Mid Cap Fund"
$ws.Cells.Item(15,4).Value = 1.49927
$ws.Cells.Item(15,5).Value = 1.496639
$ws.Cells.Item(15,6).Value = 1.505265
$ws.Cells.Item(15,7).Value = 0.00263100000000005
$ws.Cells.Item(15,8).Value = -0.005994999999999973
$ws.Cells.Item(16,1).Value = "INE881D01027"
$ws.Cells.Item(16,2).Value = "Oracle Financial Services Software Ltd"
$ws.Cells.Item(16,3).Value = "quant Mid Cap Fund"
$ws.Cells.Item(16,4).Value = 1.441495
$ws.Cells.Item(16,5).Value = 0.451678
$ws.Cells.Item(16,6).Value = 0.459577
$ws.Cells.Item(16,7).Value = 0.989817
$ws.Cells.Item(16,8).Value = 0.981918
$ws.Cells.Item(17,1).Value = "INE298J01013"
$ws.Cells.Item(17,2).Value = "Nippon Life India Asset Management Ltd"
$ws.Cells.Item(17,3).Value = "quant Mid Cap Fund"
$ws.Cells.Item(17,4).Value = 1.333105
$ws.Cells.Item(17,5).Value = 0
$ws.Cells.Item(17,6).Value = 0
$ws.Cells.Item(17,7).Value = 1.333105
$ws.Cells.Item(17,8).Value = 1.333105
$ws.Cells.Item(18,1).Value = "INE776C01039"
$ws.Cells.Item(18,2).Value = "GMR Airports Limited"
$ws.Cells.Item(18,3).Value = "quant Mid Cap Fund"
$ws.Cells.Item(18,4).Value = 1.181887
$ws.Cells.Item(18,5).Value = 2.43113
$ws.Cells.Item(18,6).Value = 5.886702
$ws.Cells.Item(18,7).Value = -1.249243
$ws.Cells.Item(18,8).Value = -4.704815
$ws.Cells.Item(19,1).Value = "INE699H01024"
$ws.Cells.Item(19,2).Value = "Adani Wilmar Limited"
$ws.Cells.Item(19,3).Value = "quant Mid Cap Fund"
$ws.Cells.Item(19,4).Value = 1.150061
$ws.Cells.Item(19,5).Value = 1.152939
$ws.Cells.Item(19,6).Value = 1.214043
$ws.Cells.Item(19,7).Value = -0.002877999999999936
$ws.Cells.Item(19,8).Value = -0.06398199999999998
$ws.Cells.Item(20,1).Value = "INE0J1Y01017"
$ws.Cells.Item(20,2).Value = "Life Insurance Corporation Of India"
$ws.Cells.Item(20,3).Value = "quant Mid Cap Fund"
$ws.Cells.Item(20,4).Value = 1.118878
$ws.Cells.Item(20,5).Value = 1.048692
$ws.Cells.Item(20,6).Value = 1.058558
$ws.Cells.Item(20,7).Value = 0.07018600000000008
$ws.Cells.Item(20,8).Value = 0.06031999999999993
$ws.Cells.Item(21,1).Value = "INE584A01023"
$ws.Cells.Item(21,2).Value = "NMDC Ltd"
$ws.Cells.Item(21,3).Value = "quant Mid Cap Fund"
$ws.Cells.Item(21,4).Value = 0.717933
$ws.Cells.Item(21,5).Value = 0
$ws.Cells.Item(21,6).Value = 0
$ws.Cells.Item(21,7).Value = 0.717933
$ws.Cells.Item(21,8).Value = 0.717933
$ws.Cells.Item(22,1).Value = "INE376G01013"
$ws.Cells.Item(22,2).Value = "Biocon Ltd"
$ws.Cells.Item(22,3).Value = "quant Mid Cap Fund"
$ws.Cells.Item(22,4).Value = 0
$ws.Cells.Item(22,5).Value = 0.381318
$ws.Cells.Item(22,6).Value = 0
$ws.Cells.Item(22,7).Value = -0.381318
$ws.Cells.Item(22,8).Value = 0
$ws.Cells.Item(23,1).Value = "INE245A01021"
$ws.Cells.Item(23,2).Value = "Tata Power Company Limited"
$ws.Cells.Item(23,3).Value = "quant Mid Cap Fund"
$ws.Cells.Item(23,4).Value = 0
$ws.Cells.Item(23,5).Value = 2.202735
$ws.Cells.Item(23,6).Value = 2.183924
$ws.Cells.Item(23,7).Value = -2.202735
$ws.Cells.Item(23,8).Value = -2.183924
$ws.Cells.Item(24,1).Value = "INE154A01025"
$ws.Cells.Item(24,2).Value = "ITC Limited"
$ws.Cells.Item(24,3).Value = "quant Mid Cap Fund"
$ws.Cells.Item(24,4).Value = 0
$ws.Cells.Item(24,5).Value = 2.500818
$ws.Cells.Item(24,6).Value = 2.420208
$ws.Cells.Item(24,7).Value = -2.500818
$ws.Cells.Item(24,8).Value = -2.420208
$ws.Cells.Item(25,1).Value = "INE115A01026"
$ws.Cells.Item(25,2).Value = "LIC Housing Finance Ltd"
$ws.Cells.Item(25,3).Value = "quant Mid Cap Fund"
$ws.Cells.Item(25,4).Value = 0
$ws.Cells.Item(25,5).Value = 1.475937
$ws.Cells.Item(25,6).Value = 1.449149
$ws.Cells.Item(25,7).Value = -1.475937
$ws.Cells.Item(25,8).Value = -1.449149
$ws.Cells.Item(26,1).Value = "INE335Y01020"
$ws.Cells.Item(26,2).Value = "Indian Railway Catering & Tourism Corp"
$ws.Cells.Item(26,3).Value = "quant Mid Cap Fund"
$ws.Cells.Item(26,4).Value = 0
$ws.Cells.Item(26,5).Value = 0
$ws.Cells.Item(26,6).Value = 1.360811
$ws.Cells.Item(26,7).Value = 0
$ws.Cells.Item(26,8).Value = -1.360811
